$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns E:H (proc_time, operating_cost, operators, operating_status) carry stale
# bestFit column-width metadata from the old layout. Wipe them out (data + widths)
# and reinsert blank replacements so the widths reset to sheet defaults, then
# repopulate - this keeps the later-columns (I..P) correctly back in place.
$ws.Range("E1:H1").EntireColumn.Delete()
$ws.Range("E1:H1").EntireColumn.Insert()

$ws.Range("E1").Value = "proc_time"
$ws.Range("F1").Value = "operating_cost"
$ws.Range("G1").Value = "operators"
$ws.Range("H1").Value = "operating_status"

$procTime = 10, 12, 14, 16, 18, 20, 22, 24, 26, 28
$opCost   = 150, 160, 170, 180, 190, 200, 210, 220, 230, 240
$operators = 1, 2, 4, 6, 1, 2, 4, 6, 1, 2

for ($i = 0; $i -lt 10; $i++) {
  $r = $i + 2
  $ws.Range("E$r").Value = $procTime[$i]
  $ws.Range("F$r").Value = $opCost[$i]
  $ws.Range("G$r").Value = $operators[$i]
  $ws.Range("H$r").Value = "available"
}

# Remove the "kwargs" column (P) entirely - header + all data
$ws.Range("P1:P11").EntireColumn.Delete()

# Column E best-fit width ~10 (matches the new narrower numeric content)
$ws.Columns.Item(5).ColumnWidth = 9.166666666666666

# Move the active selection to P5 (matches the saved view state)
$ws.Range("P5").Select() | Out-Null
